# Update "想去人数" (interested-count) values in column F for the
# "展览" (Exhibition) and "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 323
$ws1.Range("F4").Value = 256
$ws1.Range("F5").Value = 3018
$ws1.Range("F6").Value = 2027
$ws1.Range("F7").Value = 394
$ws1.Range("F8").Value = 140
$ws1.Range("F9").Value = 1128
$ws1.Range("F10").Value = 207
$ws1.Range("F11").Value = 734
$ws1.Range("F12").Value = 65

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 323
$ws4.Range("F4").Value = 256
$ws4.Range("F5").Value = 3018
$ws4.Range("F6").Value = 2027
$ws4.Range("F7").Value = 394
# F8 on this sheet corresponds to a different event (row not present on
# the 展览 sheet) and is left unchanged.
$ws4.Range("F9").Value = 140
$ws4.Range("F10").Value = 1128
$ws4.Range("F11").Value = 207
$ws4.Range("F12").Value = 734
$ws4.Range("F13").Value = 65
